# Applies the "Add additional user credentials to Postman CSV and update
# Skills.xlsx" edit to Skills.xlsx:
#   1. Fix a typo in C6 ("kỹ" -> "kỹ thuật").
#   2. Insert a new row after row 24 (inside the API section) containing a
#      new skill line ("Viết các scripts để test hồi quy & tự động"),
#      pushing every row below it down by one.
#   3. Nudge the floating picture (anchored via a one-cell anchor near the
#      bottom of the sheet) down by one row so it keeps tracking the same
#      relative position after the insert.
#   4. Re-apply the row-height bumps that show up in the saved file
#      (rows 9-12 & 14-16 -> 18pt, row 23 -> 54pt) and refresh the
#      selection/zoom the workbook was left on.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Typo fix on row 6.
$ws.Range("C6").Value = "Nắm bản chất kỹ thuật Static Testing và Dynamic Testing"

# 2) Insert the new row for the Postman regression-script bullet.
$ws.Rows("25:25").Insert()
$ws.Range("C25").Value = "Viết các scripts để test hồi quy & tự động"

# 3) Shift the floating picture down by the height of the newly inserted
#    row (default row height) so it keeps sitting just past the last used
#    row, matching where Excel leaves it after the insert.
$shp = $ws.Shapes.Item(1)
$shp.Top = $shp.Top() + $ws.Rows(25).RowHeight()

# 4) Row-height touch-ups that Excel recorded for this sheet.
$ws.Rows("9:12").RowHeight = 18
$ws.Rows("14:16").RowHeight = 18
$ws.Rows(23).RowHeight = 54

# 5) Leave the selection/zoom the way the workbook was saved.
$ws.Range("B16").Select()
$excel.ActiveWindow.Zoom = 180
